$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 corresponds to the "van Dis (2020)" metrics row being corrected
# (Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020))

$ws.Range("C3").Value = 0.8484848484848485
$ws.Range("H3").Value = 0.7649746192893401
$ws.Range("I3").Value = 0.05384771827724848
$ws.Range("J3").Value = 0.7575757575757576
$ws.Range("K3").Value = 121.3030303030303

$ws.Range("Q3").Value = 7
$ws.Range("R3").Value = 14
$ws.Range("S3").Value = 49
$ws.Range("T3").Value = 95
$ws.Range("U3").Value = 175
$ws.Range("V3").Value = 1930
$ws.Range("W3").Value = 1923
$ws.Range("X3").Value = 1888
$ws.Range("Y3").Value = 1842
$ws.Range("Z3").Value = 1762

$ws.Range("AF3").Value = 0.996386
$ws.Range("AG3").Value = 0.992772
$ws.Range("AH3").Value = 0.974703
$ws.Range("AI3").Value = 0.950955
$ws.Range("AJ3").Value = 0.909654
